$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)

# Replace the whole paragraph text with the new combined wording, then
# re-split it back into individual runs (mirrors what the author did:
# retype the "Goal of the project" paragraph with the new text about
# "KoBun").
$shape.TextFrame.DeleteText()
$shape.TextFrame.TextRange.Text = "Целью проекта ""КоBun"" является разработка игры на клетчатом поле, предоставляющей пользователю интересные и захватывающие игровые сценарии. Главной задачей игрока будет выбраться из лабиринта, избегая встреч с враждебными персонажами. "

$tr = $shape.TextFrame.TextRange

# Force "КоBun" (the misspelled/mixed-script product name) onto its own
# run, matching the proofing-flagged run in the final deck.
$nameRun = $tr.Characters(16, 5)
$nameRun.Font.Bold = $nameRun.Font.Bold

# Force the closing-quote + remainder of the sentence onto its own run
# as well, so the paragraph ends up split exactly like the source.
$tailRun = $tr.Characters(21, 215)
$tailRun.Font.Bold = $tailRun.Font.Bold
